# Update gh-pages to output generated at 456a3b4
# Bumps the "想去人数" (want-to-go count) figures in column F for a handful
# of rows on the "展览" sheet and the mirrored rows on the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 829
$wsExhibit.Range("F5").Value = 306
$wsExhibit.Range("F9").Value = 557
$wsExhibit.Range("F13").Value = 13556
$wsExhibit.Range("F17").Value = 5570
$wsExhibit.Range("F19").Value = 63

# --- Sheet "全部类型" (combined view mirroring "展览" rows) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 829
$wsAll.Range("F21").Value = 306
$wsAll.Range("F31").Value = 557
$wsAll.Range("F35").Value = 13556
$wsAll.Range("F40").Value = 5570
$wsAll.Range("F42").Value = 63
